# Apply updated sweep-measurement readings (re-run of the measurement sweep)
# to the "Measurements" sheet (row 2 + row 3), refresh the comment text, and
# recompute the per-metric Max/Min/Mean rows on the "Statistics" sheet.

$wb = $excel.ActiveWorkbook
$wsMeas  = $wb.Worksheets.Item("Measurements")
$wsStats = $wb.Worksheets.Item("Statistics")

# New values for columns A..AW (49 columns) for data rows 2 and 3.
$row2Vals = @(4.544441938400269, 13.52477216720581, 2.0, 5.0, 1.0, 0.393, -13.43059, 5.085418, 4.976383209, -47.23838425, 0.1952049732208252, 5.02750635147, -55.1757493019, -56.4646325111, 0.217111349105835, 5.057008743, -50.3416481, 0.2275807857513428, 5.10628175735, -60.6061177254, -60.3064007759, 0.1907782554626465, 6.935030698776245, 0.0, 0.0, 2.638, 4.986499786, -50.09896851, 0.3401844501495361, 5.03555059433, -55.7622961998, -55.8825850487, 0.1863663196563721, 11.75644135475159, 0.0, 0.0, 2.571, 4.984535217, -49.71195221, 0.3320817947387695, 5.04499912262, -58.7719774246, -58.9627084732, 0.1902198791503906, 21.41870021820068, 0.0, 0.0, 2.566, 41.271)
$row3Vals = @(4.544441938400269, 13.52477216720581, 2.1, 5.0, 1.0, 0.352, -13.47612, 5.036761, 4.908203125, -21.40676689, 0.1922016143798828, 4.97490406036, -55.1714296341, -56.3040742874, 0.1937520503997803, 4.908157349, -21.40340424, 0.1980037689208984, 4.98175525665, -55.0589227676, -56.2393503189, 0.185863733291626, 3.048375129699707, 0.0, 0.0, 1.218, 5.009925842, -49.20523834, 0.3543787002563477, 5.07956218719, -54.4045591354, -55.4119215012, 0.186917781829834, 11.80976676940918, 0.0, 0.0, 2.545, 5.010688782, -50.68193817, 0.3266315460205078, 5.0837225914, -60.4407744408, -60.4873366356, 0.1884052753448486, 20.66013479232788, 0.0, 0.0, 2.594, 36.657)

$numCols = $row2Vals.Length   # 49 -> columns A (1) .. AW (49)

for ($i = 0; $i -lt $numCols; $i++) {
    $col = $i + 1
    $wsMeas.Cells.Item(2, $col).Value2 = $row2Vals[$i]
    $wsMeas.Cells.Item(3, $col).Value2 = $row3Vals[$i]
}

# Updated comment text (added "and" before "external sensors").
$newComment = "The 5GNR waveform used in this test is a 10MHz UL, 30kHz SCS, 24QAM, 24RB, 0rbo configuration.`nThis test utilizes the full 5G frame.`nThe power servo is done after each DPD type to ensure accurate output power.`nThe power servo uses the NRX power meter and external sensors for power servo."
$wsMeas.Range("AX2").Value2 = $newComment
$wsMeas.Range("AX3").Value2 = $newComment

# ---------------------------------------------------------------------------
# Recompute the Statistics sheet: for every column A..AW of Measurements,
# three rows hold "<Header> - Max", "<Header> - Min", "<Header> - Mean".
# Row 2 of Statistics ("Number of Tests") is untouched (still 2 data rows).
# ---------------------------------------------------------------------------
$statRow = 3
for ($i = 0; $i -lt $numCols; $i++) {
    $col = $i + 1
    $v2 = $row2Vals[$i]
    $v3 = $row3Vals[$i]

    $maxV  = [Math]::Max($v2, $v3)
    $minV  = [Math]::Min($v2, $v3)
    $meanV = ($v2 + $v3) / 2.0

    $wsStats.Cells.Item($statRow,     2).Value2 = $maxV
    $wsStats.Cells.Item($statRow + 1, 2).Value2 = $minV
    $wsStats.Cells.Item($statRow + 2, 2).Value2 = $meanV

    $statRow = $statRow + 3
}
